# Fixed market issues, added a couple new affixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

# --- Fix market issue: M26 (skill_level_trivial for Chakra Alignment) 90 -> 180 ---
$ws.Cells.Item(26, 13).Value = 180

# --- Add new affixes as rows 38-42 ---

# Row 38: Shadow Fiend Lover
$ws.Cells.Item(38, 1).Value = "Shadow Fiend Lover"
$ws.Cells.Item(38, 2).Value = "Feel the strength of a shadow fiend crawling through the veins of your body."
$ws.Cells.Item(38, 3).Value = 0.05
$ws.Cells.Item(38, 4).Value = 0.05
$ws.Cells.Item(38, 5).Value = 0.05
$ws.Cells.Item(38, 11).Value = 15
$ws.Cells.Item(38, 12).Value = 3
$ws.Cells.Item(38, 13).Value = 8
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = 750
$ws.Cells.Item(38, 19).Value = "prefix"

# Row 39: Dark Thoughts
$ws.Cells.Item(39, 1).Value = "Dark Thoughts"
$ws.Cells.Item(39, 2).Value = "These thoughts are drifting through your head all the time. What can you do about it?"
$ws.Cells.Item(39, 3).Value = 0.08
$ws.Cells.Item(39, 4).Value = 0.08
$ws.Cells.Item(39, 5).Value = 0.08
$ws.Cells.Item(39, 11).Value = 35
$ws.Cells.Item(39, 12).Value = 15
$ws.Cells.Item(39, 13).Value = 30
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = 7000
$ws.Cells.Item(39, 19).Value = "prefix"

# Row 40: Fiathless Hate
$ws.Cells.Item(40, 1).Value = "Fiathless Hate"
$ws.Cells.Item(40, 2).Value = "There is nothing worse then the hatefilled vengance of the faithless."
$ws.Cells.Item(40, 3).Value = 0.1
$ws.Cells.Item(40, 4).Value = 0.1
$ws.Cells.Item(40, 5).Value = 0.1
$ws.Cells.Item(40, 11).Value = 38
$ws.Cells.Item(40, 12).Value = 15
$ws.Cells.Item(40, 13).Value = 30
$ws.Cells.Item(40, 18).Value = 100000
$ws.Cells.Item(40, 19).Value = "prefix"

# Row 41: Demonic Pact
$ws.Cells.Item(41, 1).Value = "Demonic Pact"
$ws.Cells.Item(41, 2).Value = "Make a demonic pact for the stats you want. Trust, it always works."
$ws.Cells.Item(41, 3).Value = 0.15
$ws.Cells.Item(41, 4).Value = 0.15
$ws.Cells.Item(41, 5).Value = 0.15
$ws.Cells.Item(41, 11).Value = 250
$ws.Cells.Item(41, 12).Value = 65
$ws.Cells.Item(41, 13).Value = 165
$ws.Cells.Item(41, 18).Value = 20300000
$ws.Cells.Item(41, 19).Value = "prefix"

# Row 42: Shadow Sands Dust
$ws.Cells.Item(42, 1).Value = "Shadow Sands Dust"
$ws.Cells.Item(42, 2).Value = "The curse of the shadow suns is one far too long to tell child, alas it will fill you with the magic you need."
$ws.Cells.Item(42, 3).Value = 0.18
$ws.Cells.Item(42, 4).Value = 0.18
$ws.Cells.Item(42, 5).Value = 0.18
$ws.Cells.Item(42, 11).Value = 400
$ws.Cells.Item(42, 12).Value = 90
$ws.Cells.Item(42, 13).Value = 200
$ws.Cells.Item(42, 18).Value = 440500000
$ws.Cells.Item(42, 19).Value = "prefix"
